# Add two new columns, I (I0) and J (IF), to the sheet.
# Row 1 holds the headers (same style as the other header cells),
# rows 2-34 hold the per-record numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1), mirror the formatting of the existing header cells ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-34 ---
$iValues = @(3,6,7,6,4,6,9,10,7,5,9,6,7,5,8,4,8,7,6,10,5,8,7,8,8,5,8,6,9,4,7,6,3)
$jValues = @(4,7,7,6,5,6,9,11,7,6,9,7,8,6,8,5,8,8,6,11,6,8,9,8,8,7,8,7,9,4,7,6,3)

for ($r = 2; $r -le 34; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
